$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column D ("META" tracking dates), mirroring columns B/F ---
$ws.Range("D2").Value = "META"

# D8:D11 get date values 45302..45305, formatted like the existing date
# column (B) -- copy that cell's number format (short date) via a
# format-only paste so the new cells share the same style record.
$ws.Range("B3").Copy()

$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("D8").Value = 45302

$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("D9").Value = 45303

$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("D10").Value = 45304

$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D11").Value = 45305

$excel.CutCopyMode = $false

# Column D sized the same as columns B/F (short-date column).
$ws.Columns("D").ColumnWidth = 9.67

# Move the active selection to A11 (matches the saved cursor position).
$ws.Range("A11").Select()
